$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 151
$wsExpo.Range("F6").Value = 9316
$wsExpo.Range("F10").Value = 1091
$wsExpo.Range("F12").Value = 75
$wsExpo.Range("F15").Value = 392
$wsExpo.Range("F17").Value = 248
$wsExpo.Range("F18").Value = 1232

# Sheet "全部类型" (All types) - same underlying events, update column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 151
$wsAll.Range("F8").Value = 9316
$wsAll.Range("F12").Value = 1091
$wsAll.Range("F14").Value = 75
$wsAll.Range("F17").Value = 392
$wsAll.Range("F19").Value = 248
$wsAll.Range("F20").Value = 1232
